$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (Hoja1 -> Semilla 4) ---
$ws.Name = "Semilla 4"

# --- Row 1: add new header cells F1 (user), G1 (password), H1 (rutaWinWap) ---
$ws.Range("F1").Value = "user"
$ws.Range("G1").Value = "password"
$ws.Range("H1").Value = "rutaWinWap"

# F1/G1 use the plain default text style already used elsewhere (style of A3/B3 etc.)
$ws.Range("A3").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# H1 needs a brand-new bold blue JetBrains Mono style. Start from the existing
# vertical-centered JetBrains Mono style (the one used by the old A8 cell) and
# then tweak Bold + Color so the resulting font/style entries are created once
# and then reused by the other rutaWinWap cells below.
$ws.Range("A8").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Font.Color = 14580521

# --- Row 2: C2 becomes a hyperlink-styled cell (like A2/B2/D2/E2); add rutaWinWap data cells ---
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("F2").Value = "CQ10960370"
$ws.Range("G2").Value = "Tigo.2022*"
$ws.Range("H2").Value = "C:\Program Files (x86)\Winwap Technologies\WinWAP for Windows 4.2\WinWAP4.exe"

$ws.Range("H1").Copy()
$ws.Range("F2:H2").PasteSpecial(-4122)

# Add the real hyperlink on C2 pointing at its own displayed URL (same pattern
# as the other URL cells in row 2).
$ws.Hyperlinks.Add($ws.Range("C2"), "http://10.69.60.77:8180/tigo-pos-web/wap/windex.wml")

# --- Rows 3-7 are unchanged ---

# --- Row 8 becomes the old "Vendedor / Cedula Cliente / MSIDN / MSI" header row ---
$ws.Range("A8").Value = "Vendedor"
$ws.Range("B8").Value = "Cedula Cliente"
$ws.Range("C8").Value = "MSIDN"
$ws.Range("D8").Value = "MSI"

$ws.Range("A1").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# --- Rows 9-13: MSISDN test data (5 rows), replacing the old rows 10-13 ---
$ws.Range("A9").Value = "10960370"
$ws.Range("B9").Value = "984108505"
$ws.Range("C9").Value = "3016875982"
$ws.Range("D9").Value = "732111198172291"

$ws.Range("A10").Value = "10960370"
$ws.Range("B10").Value = "835244140"
$ws.Range("C10").Value = "3016877591"
$ws.Range("D10").Value = "732111198172292"

$ws.Range("A11").Value = "10960370"
$ws.Range("B11").Value = "667299000"
$ws.Range("C11").Value = "3016875982"
$ws.Range("D11").Value = "732111198172291"

$ws.Range("A12").Value = "10960370"
$ws.Range("B12").Value = "835244140"
$ws.Range("C12").Value = "3016877411"
$ws.Range("D12").Value = "732111198172294"

$ws.Range("A13").Value = "10960370"
$ws.Range("B13").Value = "311615530"
$ws.Range("C13").Value = "3016876876"
$ws.Range("D13").Value = "732111198172293"

# Normalize the style of all of rows 9-13 to the plain text style already used
# by C11/D11 etc. in the original sheet.
$ws.Range("C11").Copy()
$ws.Range("A9:D13").PasteSpecial(-4122)

# --- Selection matches the author's final cursor position ---
$ws.Range("C13").Select()
